# Remove the "Catégories" row from the Non-Conformité sheet so it no longer
# prints on the NC form. This deletes entire row 21 (label "Catégories" plus
# its input cell) and shifts all rows below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NC")

$ws.Rows.Item(21).Delete()
